$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 64498
$ws.Range("I11").Value = 64498
$ws.Range("K11").Value = 64498
$ws.Range("M11").Value = -64358
$ws.Range("H18").Value = 1735.6
$ws.Range("I18").Value = 1669.5
$ws.Range("J18").Value = 2000
$ws.Range("K18").Value = 1669.5
$ws.Range("L18").Value = 2000
$ws.Range("M18").Value = -1385.5
$ws.Range("N18").Value = -2568
$ws.Range("H46").Value = 4924.875
$ws.Range("I46").Value = 2799.6667
$ws.Range("K46").Value = 8399.000100000001
$ws.Range("M46").Value = -8280.000100000001
$ws.Range("H60").Value = 4924.875
$ws.Range("I60").Value = 2799.6667
$ws.Range("K60").Value = 8399.000100000001
$ws.Range("M60").Value = -7915.000100000001
$ws.Range("H70").Value = 4600
$ws.Range("I70").Value = 3250
$ws.Range("J70").Value = 5090.909
$ws.Range("K70").Value = 9750
$ws.Range("L70").Value = 15272.727
$ws.Range("M70").Value = -9480
$ws.Range("N70").Value = -15812.727
$ws.Range("H73").Value = 4600
$ws.Range("I73").Value = 3250
$ws.Range("J73").Value = 5090.909
$ws.Range("K73").Value = 9750
$ws.Range("L73").Value = 15272.727
$ws.Range("M73").Value = -8814
$ws.Range("N73").Value = -17144.727
$ws.Range("H80").Value = 2418
$ws.Range("I80").Value = 1287.8667
$ws.Range("J80").Value = 4301.5557
$ws.Range("K80").Value = 3863.6001
$ws.Range("L80").Value = 12904.6671
$ws.Range("M80").Value = -2865.6001
$ws.Range("N80").Value = -14900.6671
$ws.Range("H83").Value = 2418
$ws.Range("I83").Value = 1287.8667
$ws.Range("J83").Value = 4301.5557
$ws.Range("K83").Value = 11590.8003
$ws.Range("L83").Value = 38714.0013
$ws.Range("M83").Value = -6598.800300000001
$ws.Range("N83").Value = -48698.0013
$ws.Range("H86").Value = 7233
$ws.Range("I86").Value = 6792.533
$ws.Range("J86").Value = 8176.857
$ws.Range("K86").Value = 6792.533
$ws.Range("L86").Value = 8176.857
$ws.Range("M86").Value = -5669.533
$ws.Range("N86").Value = -10422.857
$ws.Range("H89").Value = 7233
$ws.Range("I89").Value = 6792.533
$ws.Range("J89").Value = 8176.857
$ws.Range("K89").Value = 33962.665
$ws.Range("L89").Value = 40884.285
$ws.Range("M89").Value = -28346.665
$ws.Range("N89").Value = -52116.285
$ws.Range("H111").Value = 849.75
$ws.Range("I111").Value = 800
$ws.Range("K111").Value = 2400
$ws.Range("M111").Value = 667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1457.9166
$ws.Range("I2").Value = 1250.125
$ws.Range("K2").Value = 1250.125
$ws.Range("M2").Value = -1137.125
$ws.Range("H74").Value = 694799.9
$ws.Range("I74").Value = 752653
$ws.Range("K74").Value = 752653
$ws.Range("M74").Value = -751779
$ws.Range("H77").Value = 694799.9
$ws.Range("I77").Value = 752653
$ws.Range("K77").Value = 3763265
$ws.Range("M77").Value = -3758897
$ws.Range("H116").Value = 1457.9166
$ws.Range("I116").Value = 1250.125
$ws.Range("K116").Value = 1250.125
$ws.Range("M116").Value = 1043.875
$ws.Range("H125").Value = 18764
$ws.Range("J125").Value = 18764
$ws.Range("L125").Value = 18764
$ws.Range("N125").Value = -28604

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1457.9166
$ws.Range("I3").Value = 1250.125
$ws.Range("K3").Value = 1250.125
$ws.Range("M3").Value = -1136.125
$ws.Range("H22").Value = 392.25
$ws.Range("I22").Value = 218.5
$ws.Range("J22").Value = 566
$ws.Range("K22").Value = 218.5
$ws.Range("L22").Value = 566
$ws.Range("M22").Value = -45.5
$ws.Range("N22").Value = -912
$ws.Range("H86").Value = 3781.1765
$ws.Range("I86").Value = 3393.875
$ws.Range("J86").Value = 4125.4443
$ws.Range("K86").Value = 3393.875
$ws.Range("L86").Value = 4125.4443
$ws.Range("M86").Value = -2270.875
$ws.Range("N86").Value = -6371.4443
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").Value = $null
$ws.Range("H89").Value = 3781.1765
$ws.Range("I89").Value = 3393.875
$ws.Range("J89").Value = 4125.4443
$ws.Range("K89").Value = 16969.375
$ws.Range("L89").Value = 20627.2215
$ws.Range("M89").Value = -11353.375
$ws.Range("N89").Value = -31859.2215
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3577451.5
$ws.Range("I31").Value = 1404643.6
$ws.Range("K31").Value = 1404643.6
$ws.Range("M31").Value = -1404348.6
$ws.Range("H34").Value = 3577451.5
$ws.Range("I34").Value = 1404643.6
$ws.Range("K34").Value = 1404643.6
$ws.Range("M34").Value = -1404441.6
$ws.Range("H94").Value = 14448.75
$ws.Range("J94").Value = 2227.1428
$ws.Range("L94").Value = 2227.1428
$ws.Range("N94").Value = -3129.1428
$ws.Range("H99").Value = 20195.637
$ws.Range("I99").Value = 26731.5
$ws.Range("K99").Value = 26731.5
$ws.Range("M99").Value = -25233.5
$ws.Range("H107").Value = 531.069
$ws.Range("J107").Value = 788
$ws.Range("L107").Value = 788
$ws.Range("N107").Value = -4628
$ws.Range("H112").Value = 35000
$ws.Range("I112").Value = 35000
$ws.Range("K112").Value = 35000
$ws.Range("M112").Value = -33523
$ws.Range("H126").Value = 20195.637
$ws.Range("I126").Value = 26731.5
$ws.Range("K126").Value = 80194.5
$ws.Range("M126").Value = -77724.5
$ws.Range("H132").Value = 3448.0908
$ws.Range("I132").Value = 3168.2856
$ws.Range("K132").Value = 9504.856800000001
$ws.Range("M132").Value = -6974.856800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 32.466667
$ws.Range("I12").Value = 10.333333
$ws.Range("K12").Value = 30.999999
$ws.Range("M12").Value = 142.000001
$ws.Range("H37").Value = 46000
$ws.Range("J37").Value = 46000
$ws.Range("L37").Value = 138000
$ws.Range("N37").Value = -138224
$ws.Range("H75").Value = 2272.8
$ws.Range("I75").Value = 2200
$ws.Range("J75").Value = 2321.3333
$ws.Range("K75").Value = 6600
$ws.Range("L75").Value = 6963.999899999999
$ws.Range("M75").Value = -5602
$ws.Range("N75").Value = -8959.999899999999
$ws.Range("H78").Value = 2272.8
$ws.Range("I78").Value = 2200
$ws.Range("J78").Value = 2321.3333
$ws.Range("K78").Value = 19800
$ws.Range("L78").Value = 20891.9997
$ws.Range("M78").Value = -14808
$ws.Range("N78").Value = -30875.9997
$ws.Range("H107").Value = 3524.889
$ws.Range("I107").Value = 5000
$ws.Range("J107").Value = 3438.1177
$ws.Range("K107").Value = 15000
$ws.Range("L107").Value = 10314.3531
$ws.Range("M107").Value = -13080
$ws.Range("N107").Value = -14154.3531
$ws.Range("H129").Value = 772367.75
$ws.Range("I129").Value = 2002805.8
$ws.Range("K129").Value = 6008417.4
$ws.Range("M129").Value = -6003417.4
$ws.Range("H137").Value = 8520.117
$ws.Range("J137").Value = 10743.917
$ws.Range("L137").Value = 32231.751
$ws.Range("N137").Value = -42431.751

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2000
$ws.Range("J80").Value = 2000
$ws.Range("L80").Value = 2000
$ws.Range("N80").Value = -3996
$ws.Range("H83").Value = 2000
$ws.Range("J83").Value = 2000
$ws.Range("L83").Value = 10000
$ws.Range("N83").Value = -19984
$ws.Range("H97").Value = 2213.9167
$ws.Range("I97").Value = 1912.8334
$ws.Range("J97").Value = 2515
$ws.Range("K97").Value = 1912.8334
$ws.Range("L97").Value = 2515
$ws.Range("M97").Value = -1416.8334
$ws.Range("N97").Value = -3507

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2359.4546
$ws.Range("I68").Value = 2467.8235
$ws.Range("K68").Value = 2467.8235
$ws.Range("M68").Value = -1718.8235
$ws.Range("H71").Value = 2359.4546
$ws.Range("I71").Value = 2467.8235
$ws.Range("K71").Value = 12339.1175
$ws.Range("M71").Value = -8595.1175
$ws.Range("H137").Value = 74700
$ws.Range("J137").Value = 74700
$ws.Range("L137").Value = 74700
$ws.Range("N137").Value = -84900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 352.8
$ws.Range("I100").Value = 336.44446
$ws.Range("K100").Value = 672.88892
$ws.Range("M100").Value = -131.88892
$ws.Range("H112").Value = 16000
$ws.Range("J112").Value = 16000
$ws.Range("L112").Value = 16000
$ws.Range("N112").Value = -18954
$ws.Range("H126").Value = 1416.5
$ws.Range("I126").Value = 999.6667
$ws.Range("K126").Value = 2999.0001
$ws.Range("M126").Value = -529.0001000000002
